## comment from python script
$wb = $excel.ActiveWorkbook

# --- Astronauta sheet: add new attendance/check marks (value 1) ---
$wsAstro = $wb.Worksheets.Item("Astronauta")
$wsAstro.Range("F3").Value = 1
$wsAstro.Range("F9").Value = 1
$wsAstro.Range("E11").Value = 1
$wsAstro.Range("F12").Value = 1
$wsAstro.Range("G14").Value = 1
$wsAstro.Range("F15").Value = 1
$wsAstro.Range("E17").Value = 1

# --- Senador sheet: add new marks (value 0) and clear the stray "Ñ" text cell ---
$wsSenador = $wb.Worksheets.Item("Senador")
$wsSenador.Range("F4").Value = 0
$wsSenador.Range("G4").Value = 0
$wsSenador.Range("H4").Value = 0
$wsSenador.Range("F11").Value = 0
$wsSenador.Range("F14").ClearContents()

# --- Update selections on Senador, Ninja before finally landing on Astronauta ---
$wsSenador.Activate() | Out-Null
$wsSenador.Range("G11").Select() | Out-Null

$wsNinja = $wb.Worksheets.Item("Ninja")
$wsNinja.Activate() | Out-Null
$wsNinja.Range("B22").Select() | Out-Null

# --- Astronauta becomes the active tab with selection F11 ---
$wsAstro.Activate() | Out-Null
$wsAstro.Range("F11").Select() | Out-Null
